# "Added Form with mock data" - populate the Receipts sheet with mock
# contribution data and split the Name column into Last Name / First Name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Receipts")

# --- Non-name header labels (reuse existing shared strings; column order
# already matches the final layout apart from the Name -> Last/First split
# handled further below) ---
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "City"
$ws.Range("E1").Value = "State"
$ws.Range("F1").Value = "Zip"
$ws.Range("G1").Value = "Contribution Source"
$ws.Range("H1").Value = "Contribution Type"
$ws.Range("I1").Value = "Occupation/Employer"
$ws.Range("J1").Value = "Date"
$ws.Range("K1").Value = "Amount"
$ws.Range("L1").Value = "Fees"

# --- Row 2: PT-PAC of MO monetary contribution ---
$ws.Range("A2").Value = "PT-PAC of MO"
$ws.Range("C2").Value = "205 E Capitol"
$ws.Range("D2").Value = "Jefferson City"
$ws.Range("E2").Value = "MO"
$ws.Range("F2").Value = 65101
$ws.Range("G2").Value = "PAC"
$ws.Range("H2").Value = "Monetary"
$ws.Range("K2").Value = 350

# --- Row 3: PT-PAC of MO monetary contribution (earlier date) ---
$ws.Range("A3").Value = "PT-PAC of MO"
$ws.Range("C3").Value = "205 E Capitol"
$ws.Range("D3").Value = "Jefferson City"
$ws.Range("E3").Value = "MO"
$ws.Range("F3").Value = 65101
$ws.Range("G3").Value = "PAC"
$ws.Range("H3").Value = "Monetary"
$ws.Range("K3").Value = 250

# --- Row 4: A Better Missouri monetary contribution ---
$ws.Range("A4").Value = "A Better Missouri"
$ws.Range("C4").Value = "205 E Capitol"
$ws.Range("D4").Value = "Jefferson City"
$ws.Range("E4").Value = "MO"
$ws.Range("F4").Value = 65101
$ws.Range("G4").Value = "PAC"
$ws.Range("H4").Value = "Monetary"
$ws.Range("K4").Value = 500

# --- Row 5: individual in-kind contribution (name filled in further below) ---
$ws.Range("C5").Value = "48 Willmore Rd"
$ws.Range("D5").Value = "St Louis"
$ws.Range("E5").Value = "MO"
$ws.Range("F5").Value = 63109
$ws.Range("G5").Value = "Individual"
$ws.Range("H5").Value = "In-Kind"
$ws.Range("I5").Value = "City of St Louis"
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 2.5

# --- Split the single "Name" column into "Last Name" / "First Name" ---
$ws.Range("A1").Value = "Last Name"
$ws.Range("B1").Value = "First Name"
$ws.Range("A5").Value = "Bushmeyer"
$ws.Range("B5").Value = "Ed"

# --- Contribution dates (built-in date format, shared across the column) ---
$ws.Range("J2").NumberFormat = "mm-dd-yy"
$ws.Range("J2").Copy()
$ws.Range("J3:J5").PasteSpecial(-4122)
$ws.Range("J2").Value = [DateTime]"2019-08-20"
$ws.Range("J3").Value = [DateTime]"2019-07-20"
$ws.Range("J4").Value = [DateTime]"2019-08-20"
$ws.Range("J5").Value = [DateTime]"2019-06-21"

# --- Selection moves to the first empty row below the new data ---
$ws.Range("A6").Select()
